$d = $word.ActiveDocument

# 1. Replace the long "expected project outcomes" paragraph text.
$old1 = "The expected project outcomes include the development of a highly refined and effective model designed specifically for business use. This model will play a crucial role in combating spam within email systems. By efficiently filtering out spam, the model will not only eliminate irrelevant or nuisance communications but also significantly reduce the risks associated with phishing attempts. These phishing attempts often target sensitive company data and finances. As such, the implementation of this model will enhance the overall security and efficiency of the business's email communications."
$new1 = "The primary goal of the project is to develop a machine learning model that specializes in identifying and classifying email content as SPAM or non-SPAM (HAM). This model, trained on the 2007 TREC Public Spam Corpus and Enron emails dataset, aims to achieve an accuracy rate of 95% or higher in its classification tasks, thereby significantly enhancing email security and efficiency within business contexts."

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

Write-Output "done"
